# Add a new weekly ranking sheet "2025-08-18" after the last existing sheet,
# positioned at the end, matching the layout/style of the prior week sheets.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-08-18"

# Header row (rank/title/author/latest_episode), bold + centered + boxed,
# matching the header style used on every other weekly sheet.
$ws.Range("A1").Value = 'rank'
$ws.Range("B1").Value = 'title'
$ws.Range("C1").Value = 'author'
$ws.Range("D1").Value = 'latest_episode'
$headerRng = $ws.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Ranking data rows 2-51 (rank, title, author, latest_episode)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = '宇崎ちゃんは遊びたい！'
$ws.Range("C2").Value = '丈(著者)'
$ws.Range("D2").Value = '第126話'
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = '魔術師クノンは見えている'
$ws.Range("C3").Value = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$ws.Range("D3").Value = '第39話②'
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～'
$ws.Range("C4").Value = 'zunta(作画) はらわたさいぞう(原作)'
$ws.Range("D4").Value = '第31話：完全なる死角②'
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$ws.Range("C5").Value = '作画：マエD 原作：新人'
$ws.Range("D5").Value = '第5話(4)'
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = '不徳のギルド'
$ws.Range("C6").Value = '河添太一'
$ws.Range("D6").Value = '第９７話：立派に育った所'
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = '蜘蛛ですが、なにか？'
$ws.Range("C7").Value = 'かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)'
$ws.Range("D7").Value = '第76話その1'
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'みつばものがたり 呪いの少女と死の輪舞《ロンド》'
$ws.Range("C8").Value = '堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)'
$ws.Range("D8").Value = '第11話：勝利の美酒'
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = '最強の少年聖騎士、転生者を狩る'
$ws.Range("C9").Value = '作画：御塩 原作：宇奈木ユラ'
$ws.Range("D9").Value = '第7話(1)'
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$ws.Range("C10").Value = '光永康則'
$ws.Range("D10").Value = '第６８話『施錠停止』②'
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = '最強勇者パーティーは愛が知りたい'
$ws.Range("C11").Value = '山田肌襦袢'
$ws.Range("D11").Value = '第29話「きみがきみであればいい」'
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$ws.Range("C12").Value = 'マツモトケンゴ'
$ws.Range("D12").Value = '第６２話　尋問の戦いが始まった（２）'
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = '生徒会にも穴はある！'
$ws.Range("C13").Value = 'むちまろ'
$ws.Range("D13").Value = '第133話'+"`t"+'みんなの将来！'
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'クラスで２番目に可愛い女の子と友だちになった'
$ws.Range("C14").Value = '尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)'
$ws.Range("D14").Value = '第34話②'
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'いとこのこ'
$ws.Range("C15").Value = 'いぬちく(著者)'
$ws.Range("D15").Value = '第37話'
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～'
$ws.Range("C16").Value = '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)'
$ws.Range("D16").Value = '第10話-1'
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = '実は俺、最強でした？'
$ws.Range("C17").Value = '原作：澄守 彩 漫画：高橋 愛'
$ws.Range("D17").Value = '第122話　因縁の再会'
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = '帰ってください！ 阿久津さん'
$ws.Range("C18").Value = '長岡太一(著者)'
$ws.Range("D18").Value = '第194話'
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = '辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？'
$ws.Range("C19").Value = 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)'
$ws.Range("D19").Value = '第10話'
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$ws.Range("C20").Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$ws.Range("D20").Value = '第５１話　英雄を倒す器用貧乏（３）'
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 'えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話'
$ws.Range("C21").Value = '佐藤36(作画) 煮豆シューター(原作)'
$ws.Range("D21").Value = '第4話前半'
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = '異世界魔王と召喚少女の奴隷魔術'
$ws.Range("C22").Value = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$ws.Range("D22").Value = '第127話　戦争を終わらせてみるⅢ（中編）'
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = '転生貴族の異世界冒険録 ～自重を知らない神々の使徒～'
$ws.Range("C23").Value = '夜州 nini 藻'
$ws.Range("D23").Value = '第69話'
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = '望まぬ不死の冒険者'
$ws.Range("C24").Value = '中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）'
$ws.Range("D24").Value = '第59話　ヴィステルヤ（後編）'
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = '氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic'
$ws.Range("C25").Value = '漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca'
$ws.Range("D25").Value = '第9話'
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 'リビルドワールド'
$ws.Range("C26").Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$ws.Range("D26").Value = '第72話②'
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Range("C27").Value = '内々けやき あし カオミン'
$ws.Range("D27").Value = '第137話 よくわからないけれど脱出するみたいです（２）'
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = '路地裏で拾った女の子がバッドエンド後の乙女ゲームのヒロインだった件'
$ws.Range("C28").Value = 'カボチャマスク(原作) 樋乃えなが(作画) へいろー(キャラクター原案)'
$ws.Range("D28").Value = '第1話'
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$ws.Range("C29").Value = '村上よしゆき 茨木野 あるてら'
$ws.Range("D29").Value = '第４１話　勇者、人魚王国を救い、歓迎される。あと、六邪神将が、全員来る（３）'
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = '願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜'
$ws.Range("C30").Value = 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)'
$ws.Range("D30").Value = '第5話-2：魔法のお勉強'
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = '小林さんちのメイドラゴン'
$ws.Range("C31").Value = 'クール教信者'
$ws.Range("D31").Value = '第148話'
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = '姫様“拷問”の時間です'
$ws.Range("C32").Value = '原作:春原ロビンソン　漫画:ひらけい'
$ws.Range("D32").Value = '拷問147'
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'サーシャちゃんとクラスメイトオタクくん'
$ws.Range("C33").Value = 'はぐはぐ(著者)'
$ws.Range("D33").Value = '第83話'
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'くらいあの子としたいこと'
$ws.Range("C34").Value = '碇マナツ(著者)'
$ws.Range("D34").Value = '第82話'
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$ws.Range("C35").Value = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$ws.Range("D35").Value = '第81話その3'
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'ダメ人間の愛しかた'
$ws.Range("C36").Value = '岩葉(著者)'
$ws.Range("D36").Value = '第19話前編　ダメ人間と新生活の彼女'
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$ws.Range("C37").Value = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$ws.Range("D37").Value = 'コミックス４巻ついに発売!!!'
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$ws.Range("C38").Value = '六志麻あさ 業務用餅 kisui'
$ws.Range("D38").Value = '第７０話'
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 'うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。'
$ws.Range("C39").Value = '三上こた こばやし少女 寝子空兄 ゆがー'
$ws.Range("D39").Value = '第1話　災禍の悪夢'
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 'バーサス'
$ws.Range("C40").Value = '原作：ONE 漫画：あずま京太郎 構成：bose'
$ws.Range("D40").Value = '第27話 幸せの在り処（1）'
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = '理想のヒモ生活'
$ws.Range("C41").Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Range("D41").Value = '第86話　その3'
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = '残念女幹部ブラックジェネラルさん'
$ws.Range("C42").Value = 'jin(著者)'
$ws.Range("D42").Value = '第二百話'
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = '北斗の拳 世紀末ドラマ撮影伝'
$ws.Range("C43").Value = '原案/武論尊・原哲夫 漫画/倉尾宏'
$ws.Range("D43").Value = '第74話 サウザーの右足危機一髪‼︎'
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = '俺は星間国家の悪徳領主！'
$ws.Range("C44").Value = '灘島かい（漫画） 三嶋与夢（原作） 高峰ナダレ（キャラクター原案）'
$ws.Range("D44").Value = '第40話　詳しく聞かせろ'
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 'ギャルとダンジョンと周回遅れの探索英雄譚'
$ws.Range("C45").Value = '漫画家： 水田ケンジ 原作：榊一郎 キャラクター原案：黒獅子'
$ws.Range("D45").Value = '第2話'
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = '異世界でも無難に生きたい症候群'
$ws.Range("C46").Value = '原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう'
$ws.Range("D46").Value = '第31話①'
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = '最弱貴族に転生したので悪役たちを集めてみた'
$ws.Range("C47").Value = '空野進 sorani ファルまろ'
$ws.Range("D47").Value = '第11話　最弱貴族、背中を流してもらう（２）'
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = '魔都精兵のスレイブ'
$ws.Range("C48").Value = '原作:タカヒロ　漫画:竹村洋平'
$ws.Range("D48").Value = '第159話　神奴隷'
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = '魔眼と弾丸を使って異世界をぶち抜く！'
$ws.Range("C49").Value = '漫画：瀬菜モナコ 原作：かたなかじ キャラクター原案：赤井てら'
$ws.Range("D49").Value = '第20話'
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 'スキル【万物支配】に目覚めたおっさんは、ダンジョンで生計を立てることにしました～無職から始める支配者無双～'
$ws.Range("C50").Value = '岸本和葉 原田 臙 シミズヒロノリ 吉武'
$ws.Range("D50").Value = '第5話　一撃(前編)'
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = '転生したらスライムだった件 異聞 ～魔国暮らしのトリニティ～'
$ws.Range("C51").Value = '伏瀬 戸野タエ みっつばー'
$ws.Range("D51").Value = '第108話　開国祭開幕［その3］'
